$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.265.86'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '1.564.23'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.20'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -0.50%  '
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").Value = '1.787.19'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '1.565.79'
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").Value = '27.247.41'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '218.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.13%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0704'
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0470'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").Value = '1.455.44'
$ws.Range("E34").Value = '  +1.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.37%  '
$ws.Range("E36").Value = '  +1.15%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.542'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.815'
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.980'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("E46").Value = '  +0.70%  '
$ws.Range("D47").Value = '1.700.59'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0524'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0947'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.48%  '
